$d = $word.ActiveDocument

# Locate the paragraph that contains the closing sentence of the existing
# section ("...atributos:") so we can insert the new paragraphs right after it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*respectivos atributos:*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate target paragraph ending in 'atributos:'"
}

# Insert at the end of the paragraph's range (this is *after* the paragraph
# mark, so the existing paragraph's text/run is left completely untouched).
$insPoint = $target.Range.End
$insRng = $d.Range($insPoint, $insPoint)

$rPrArial = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$pPrArial = "<w:pPr>$rPrArial</w:pPr>"
$pPrArialUnderline = '<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr>'

function Word-Run([string]$text, [bool]$preserve) {
    if ($preserve) {
        return "<w:r>$rPrArial<w:t xml:space=`"preserve`">$text</w:t></w:r>"
    } else {
        return "<w:r>$rPrArial<w:t>$text</w:t></w:r>"
    }
}

function Spell-Run([string]$text, [bool]$preserve) {
    return '<w:proofErr w:type="spellStart"/>' + (Word-Run $text $preserve) + '<w:proofErr w:type="spellEnd"/>'
}

# Paragraph 2: the long paragraph describing the tables/inserts created.
$p2runs = ""
$p2runs += Word-Run "1 creé la tabla " $true
$p2runs += Spell-Run "grado_dificultad" $false
$p2runs += Word-Run ". E " $true
$p2runs += Spell-Run "inserts" $false
$p2runs += Word-Run ", creé tabla puntuación e " $true
$p2runs += Spell-Run "inserts" $false
$p2runs += Word-Run ", creé " $true
$p2runs += Spell-Run "type" $false
$p2runs += Word-Run " para nombre, creé tabla persona que es padre y tabla líder, posteriormente hice " $true
$p2runs += Spell-Run "inserts" $false
$p2runs += Word-Run " a líder, vemos que se guardan tmb en persona. " $true

$bodyXml = ""
$bodyXml += "<w:p>$pPrArial</w:p>"
$bodyXml += "<w:p>$pPrArial$p2runs</w:p>"
$bodyXml += "<w:p>$pPrArial</w:p>"
$bodyXml += "<w:p>$pPrArialUnderline</w:p>"

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insRng.InsertXML($xml)
